$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 229, shifting existing rows 229:308 down to 230:309.
$ws.Rows(229).Insert()

# Populate the newly inserted row 229 with the new data record.
$ws.Range("A229").Value = 3
$ws.Range("B229").Value = "Femacal de La Calera"
$ws.Range("C229").Value = "Coquimbo"
$ws.Range("D229").Value = 44627
$ws.Range("D229").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E229").Value = 5
$ws.Range("F229").Value = 100112043
$ws.Range("G229").Value = "Pepino ensalada"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 110
$ws.Range("K229").Value = 15000
$ws.Range("L229").Value = 15500
$ws.Range("M229").Value = 15136
$ws.Range("N229").Value = "$/caja 70 unidades"
$ws.Range("O229").Value = "Provincia de Quillota"
$ws.Range("P229").Value = 216
$ws.Range("Q229").Value = 70
$ws.Range("R229").Value = "Hortaliza"
